$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.62233521090399, 2.522301455612483, 0.40131844, 0.009669999999999956, 23, 0),
    @(-0.1439873785766933, -1.471887876755644, 0.05802200000000002, 0.03868100000000002, 2, 1),
    @(-0.1944567417772962, -1.451368869383521, 0.062857, 0.03868100000000002, 2, 2),
    @(-0.1878290432683844, -1.455779041389828, 0.05802200000000002, 0.043516, 2, 3),
    @(0.03162457190715398, -1.450085130149375, 0.05318699999999998, 0.043516, 2, 4),
    @(-0.00173206329599508, -1.458793951136022, 0.05318699999999998, 0.03868100000000002, 2, 5),
    @(0.03567702657627076, -1.462300797227134, 0.05318699999999998, 0.03868100000000002, 2, 6),
    @(1.029235511835111, -0.3981036020233515, 0.05318699999999998, 0.778461, 6, 7),
    @(1.011012561768807, -0.373192743430923, 0.05318699999999998, 0.7736259999999999, 5, 8),
    @(1.012426788620537, -0.4290224664953883, 0.05318699999999998, 0.778461, 4, 9),
    @(0.6925130488285822, -0.5266977223525431, 0.04835200000000001, 0.04351699999999997, 5, 10),
    @(0.7774594708806193, -0.4060413732233763, 0.04835200000000001, 0.04351699999999997, 5, 11),
    @(0.7211200139768806, -0.4353120919130355, 0.04351700000000003, 0.04351699999999997, 6, 12),
    @(0.2208470512434584, -1.21268651768271, 0.04351700000000003, 0.75912, 5, 13),
    @(0.1942641042617906, -1.333907577353858, 0.04835200000000001, 0.08219700000000002, 5, 14),
    @(0.2666114568683438, -1.192584240931173, 0.04351700000000003, 0.7542849999999999, 5, 15),
    @(0.05044819977532972, -1.388127460627619, 0.773625, 0.03384599999999999, 6, 16),
    @(0.1055700686335244, -1.390640771950284, 0.06769199999999997, 0.03384599999999999, 5, 17),
    @(0.1991841031719053, -1.160178407400638, 0.06769199999999997, 0.03868100000000002, 6, 18),
    @(0.5455128807247918, -0.6817981402625692, 0.78813, 0.03384599999999999, 6, 19),
    @(0.6794565355798995, -0.7003349853864269, 0.06769200000000003, 0.043516, 5, 20),
    @(0.5281096416115046, -0.8283812768528249, 0.06769199999999997, 0.03384599999999999, 7, 21),
    @(0.7513295979020571, -0.8958518757386327, 0.865493, 0.043516, 7, 22),
    @(0.7063829281580859, -1.00022629969505, 0.06769199999999997, 0.03868100000000002, 7, 23),
    @(0.1363098416758242, -1.28884224166924, 0.07252700000000001, 0.043516, 4, 24),
    @(-0.0008673217820262162, -1.312388979180539, 0.062857, 0.03868100000000002, 5, 25),
    @(0.05228232546424906, -1.39812640837186, 0.06769199999999997, 0.03868100000000002, 3, 26),
    @(0.02658874061189204, -1.370304732648153, 0.06769199999999997, 0.043516, 5, 27),
    @(-0.0009109289802389225, -1.3161154272228, 0.062857, 0.03868100000000002, 5, 28),
    @(-0.04036099824645148, -1.281405104458613, 0.062857, 0.03868100000000002, 4, 29),
    @(1.025137790954911, -0.5870453731317928, 0.391647936, 0.01933999999999997, 8, 30),
    @(-0.2412115450904691, -1.341483288709219, 0.04835099999999998, 0.01450499999999999, 1, 31),
    @(-0.2115760351813951, -1.315507902970704, 0.831647, 0.01450499999999999, 1, 32),
    @(-0.2327186457733519, -1.313434601935141, 0.04835099999999998, 0.01450499999999999, 1, 33),
    @(-0.1908920539494794, -1.34118999317548, 0.043516, 0.01450499999999999, 2, 34),
    @(-0.3705283647899347, -1.278218847787483, 0.043516, 0.009669999999999956, 2, 35),
    @(-0.2693381946138951, -1.29476132194428, 0.043516, 0.01450499999999999, 1, 36),
    @(-0.3094347515376059, -1.302681212480544, 0.04835099999999998, 0.01450499999999999, 2, 37),
    @(0.7767575072402623, -1.101742949659991, 0.6092299999999999, 0.01933999999999997, 7, 38),
    @(-0.1307490918608382, -1.349952946420109, 0.05318599999999996, 0.01450499999999999, 2, 39),
    @(-0.2544077109809889, -1.310803800304371, 0.04835099999999998, 0.01450499999999999, 1, 40)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
